$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3538  # H98: 3607.25 -> 3538
$ws.Cells.Item(98, 9).Value = 2304.08  # I98: 2333.4583 -> 2304.08
$ws.Cells.Item(98, 11).Value = 2304.08  # K98: 2333.4583 -> 2304.08
$ws.Cells.Item(98, 13).Value = -806.0799999999999  # M98: -835.4582999999998 -> -806.0799999999999
$ws.Cells.Item(122, 8).Value = 3538  # H122: 3607.25 -> 3538
$ws.Cells.Item(122, 9).Value = 2304.08  # I122: 2333.4583 -> 2304.08
$ws.Cells.Item(122, 11).Value = 6912.24  # K122: 7000.374899999999 -> 6912.24
$ws.Cells.Item(122, 13).Value = -4462.24  # M122: -4550.374899999999 -> -4462.24
$ws.Cells.Item(138, 8).Value = 5043.6665  # H138: 5129.2925 -> 5043.6665
$ws.Cells.Item(138, 10).Value = 5672.944  # J138: 6069.061 -> 5672.944
$ws.Cells.Item(138, 12).Value = 17018.832  # L138: 18207.183 -> 17018.832
$ws.Cells.Item(138, 14).Value = -27298.832  # N138: -28487.183 -> -27298.832

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 724.16  # H97: 750.2857 -> 724.16
$ws.Cells.Item(97, 9).Value = 713.125  # I97: 750.2857 -> 713.125
$ws.Cells.Item(97, 10).Value = 989  # J97: 0 -> 989
$ws.Cells.Item(97, 11).Value = 713.125  # K97: 750.2857 -> 713.125
$ws.Cells.Item(97, 12).Value = 989  # L97: 0 -> 989
$ws.Cells.Item(97, 13).Value = -217.125  # M97: -254.2857 -> -217.125
$ws.Cells.Item(97, 14).Value = -1981  # N97: None -> -1981
$ws.Cells.Item(122, 8).Value = 3563.6365  # H122: 3753 -> 3563.6365
$ws.Cells.Item(122, 9).Value = 2050  # I122: 1008 -> 2050
$ws.Cells.Item(122, 10).Value = 4428.5713  # J122: 5400 -> 4428.5713
$ws.Cells.Item(122, 11).Value = 6150  # K122: 3024 -> 6150
$ws.Cells.Item(122, 12).Value = 13285.7139  # L122: 16200 -> 13285.7139
$ws.Cells.Item(122, 13).Value = -3700  # M122: -574 -> -3700
$ws.Cells.Item(122, 14).Value = -18185.7139  # N122: -21100 -> -18185.7139
$ws.Cells.Item(124, 8).Value = 34990  # H124: 37000 -> 34990
$ws.Cells.Item(124, 10).Value = 34990  # J124: 37000 -> 34990
$ws.Cells.Item(124, 12).Value = 34990  # L124: 37000 -> 34990
$ws.Cells.Item(124, 14).Value = -44810  # N124: -46820 -> -44810
$ws.Cells.Item(125, 8).Value = 43736.668  # H125: 43354 -> 43736.668
$ws.Cells.Item(125, 10).Value = 43736.668  # J125: 43354 -> 43736.668
$ws.Cells.Item(125, 12).Value = 43736.668  # L125: 43354 -> 43736.668
$ws.Cells.Item(125, 14).Value = -53576.668  # N125: -53194 -> -53576.668
$ws.Cells.Item(127, 8).Value = 0  # H127: 42382.855 -> 0
$ws.Cells.Item(127, 10).Value = 0  # J127: 42382.855 -> 0
$ws.Cells.Item(127, 12).Value = 0  # L127: 42382.855 -> 0
$ws.Cells.Item(127, 14).Value = ""  # N127: -52302.855 -> None
$ws.Cells.Item(128, 8).Value = 0  # H128: 42580 -> 0
$ws.Cells.Item(128, 10).Value = 0  # J128: 42580 -> 0
$ws.Cells.Item(128, 12).Value = 0  # L128: 42580 -> 0
$ws.Cells.Item(128, 14).Value = ""  # N128: -52540 -> None
$ws.Cells.Item(130, 8).Value = 48114.5  # H130: 35667 -> 48114.5
$ws.Cells.Item(130, 10).Value = 48114.5  # J130: 35667 -> 48114.5
$ws.Cells.Item(130, 12).Value = 48114.5  # L130: 35667 -> 48114.5
$ws.Cells.Item(130, 14).Value = -58154.5  # N130: -45707 -> -58154.5
$ws.Cells.Item(131, 8).Value = 39715  # H131: 39760.832 -> 39715
$ws.Cells.Item(131, 10).Value = 39715  # J131: 39760.832 -> 39715
$ws.Cells.Item(131, 12).Value = 39715  # L131: 39760.832 -> 39715
$ws.Cells.Item(131, 14).Value = -49795  # N131: -49840.832 -> -49795
$ws.Cells.Item(134, 8).Value = 60000  # H134: 44960.832 -> 60000
$ws.Cells.Item(134, 10).Value = 60000  # J134: 44960.832 -> 60000
$ws.Cells.Item(134, 12).Value = 60000  # L134: 44960.832 -> 60000
$ws.Cells.Item(134, 14).Value = -70140  # N134: -55100.832 -> -70140

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 3601.3333  # H12: 652.5 -> 3601.3333
$ws.Cells.Item(12, 9).Value = 902.5  # I12: 652.5 -> 902.5
$ws.Cells.Item(12, 10).Value = 8999  # J12: 0 -> 8999
$ws.Cells.Item(12, 11).Value = 902.5  # K12: 652.5 -> 902.5
$ws.Cells.Item(12, 12).Value = 8999  # L12: 0 -> 8999
$ws.Cells.Item(12, 13).Value = -734.5  # M12: -484.5 -> -734.5
$ws.Cells.Item(12, 14).Value = -9335  # N12: None -> -9335
$ws.Cells.Item(99, 8).Value = 4130  # H99: 4923.846 -> 4130
$ws.Cells.Item(99, 9).Value = 2001.4286  # I99: 2577.5 -> 2001.4286
$ws.Cells.Item(99, 10).Value = 5620  # J99: 5966.6665 -> 5620
$ws.Cells.Item(99, 11).Value = 2001.4286  # K99: 2577.5 -> 2001.4286
$ws.Cells.Item(99, 12).Value = 5620  # L99: 5966.6665 -> 5620
$ws.Cells.Item(99, 13).Value = -503.4286  # M99: -1079.5 -> -503.4286
$ws.Cells.Item(99, 14).Value = -8616  # N99: -8962.666499999999 -> -8616
$ws.Cells.Item(105, 8).Value = 37039730  # H105: 41669410 -> 37039730
$ws.Cells.Item(105, 9).Value = 41669070  # I105: 47621464 -> 41669070
$ws.Cells.Item(105, 11).Value = 41669070  # K105: 47621464 -> 41669070
$ws.Cells.Item(105, 13).Value = -41667323  # M105: -47619717 -> -41667323

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2988  # H99: 3005.1562 -> 2988
$ws.Cells.Item(99, 9).Value = 1873.1904  # I99: 1877.9524 -> 1873.1904
$ws.Cells.Item(99, 10).Value = 5116.273  # J99: 5157.091 -> 5116.273
$ws.Cells.Item(99, 11).Value = 1873.1904  # K99: 1877.9524 -> 1873.1904
$ws.Cells.Item(99, 12).Value = 5116.273  # L99: 5157.091 -> 5116.273
$ws.Cells.Item(99, 13).Value = -375.1904  # M99: -379.9523999999999 -> -375.1904
$ws.Cells.Item(99, 14).Value = -8112.273  # N99: -8153.091 -> -8112.273
$ws.Cells.Item(122, 8).Value = 4260.8667  # H122: 4422.357 -> 4260.8667
$ws.Cells.Item(122, 10).Value = 4745.8887  # J122: 5089.125 -> 4745.8887
$ws.Cells.Item(122, 12).Value = 14237.6661  # L122: 15267.375 -> 14237.6661
$ws.Cells.Item(122, 14).Value = -19137.6661  # N122: -20167.375 -> -19137.6661
$ws.Cells.Item(126, 8).Value = 2988  # H126: 3005.1562 -> 2988
$ws.Cells.Item(126, 9).Value = 1873.1904  # I126: 1877.9524 -> 1873.1904
$ws.Cells.Item(126, 10).Value = 5116.273  # J126: 5157.091 -> 5116.273
$ws.Cells.Item(126, 11).Value = 5619.5712  # K126: 5633.857199999999 -> 5619.5712
$ws.Cells.Item(126, 12).Value = 15348.819  # L126: 15471.273 -> 15348.819
$ws.Cells.Item(126, 13).Value = -3149.5712  # M126: -3163.857199999999 -> -3149.5712
$ws.Cells.Item(126, 14).Value = -20288.819  # N126: -20411.273 -> -20288.819
$ws.Cells.Item(132, 8).Value = 3953.8125  # H132: 4928.7393 -> 3953.8125
$ws.Cells.Item(132, 9).Value = 3638.7222  # I132: 6356.5 -> 3638.7222
$ws.Cells.Item(132, 10).Value = 4358.9287  # J132: 4167.2666 -> 4358.9287
$ws.Cells.Item(132, 11).Value = 10916.1666  # K132: 19069.5 -> 10916.1666
$ws.Cells.Item(132, 12).Value = 13076.7861  # L132: 12501.7998 -> 13076.7861
$ws.Cells.Item(132, 13).Value = -8386.1666  # M132: -16539.5 -> -8386.1666
$ws.Cells.Item(132, 14).Value = -18136.7861  # N132: -17561.7998 -> -18136.7861

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 17865700  # H131: 15158792 -> 17865700
$ws.Cells.Item(131, 9).Value = 166736670  # I131: 50021252 -> 166736670
$ws.Cells.Item(131, 10).Value = 1183.2  # J131: 1199.3043 -> 1183.2
$ws.Cells.Item(131, 11).Value = 500210010  # K131: 150063756 -> 500210010
$ws.Cells.Item(131, 12).Value = 3549.6  # L131: 3597.9129 -> 3549.6
$ws.Cells.Item(131, 13).Value = -500204970  # M131: -150058716 -> -500204970
$ws.Cells.Item(131, 14).Value = -13629.6  # N131: -13677.9129 -> -13629.6
$ws.Cells.Item(133, 8).Value = 3434.4  # H133: 3435.2 -> 3434.4
$ws.Cells.Item(133, 9).Value = 3348.5715  # I133: 3392 -> 3348.5715
$ws.Cells.Item(133, 10).Value = 3543.6365  # J133: 3500 -> 3543.6365
$ws.Cells.Item(133, 11).Value = 10045.7145  # K133: 10176 -> 10045.7145
$ws.Cells.Item(133, 12).Value = 10630.9095  # L133: 10500 -> 10630.9095
$ws.Cells.Item(133, 13).Value = -4985.7145  # M133: -5116 -> -4985.7145
$ws.Cells.Item(133, 14).Value = -20750.9095  # N133: -20620 -> -20750.9095
$ws.Cells.Item(137, 8).Value = 9988.385  # H137: 10986.454 -> 9988.385
$ws.Cells.Item(137, 10).Value = 10580.591  # J137: 11932.056 -> 10580.591
$ws.Cells.Item(137, 12).Value = 31741.773  # L137: 35796.16800000001 -> 31741.773
$ws.Cells.Item(137, 14).Value = -41941.773  # N137: -45996.16800000001 -> -41941.773
$ws.Cells.Item(140, 8).Value = 64557.5  # H140: 51666 -> 64557.5
$ws.Cells.Item(140, 9).Value = 64557.5  # I140: 57073.332 -> 64557.5
$ws.Cells.Item(140, 10).Value = 0  # J140: 3000 -> 0
$ws.Cells.Item(140, 11).Value = 193672.5  # K140: 171219.996 -> 193672.5
$ws.Cells.Item(140, 12).Value = 0  # L140: 9000 -> 0
$ws.Cells.Item(140, 13).Value = -188492.5  # M140: -166039.996 -> -188492.5
$ws.Cells.Item(140, 14).Value = ""  # N140: -19360 -> None

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 62502250  # H80: 50002400 -> 62502250
$ws.Cells.Item(83, 8).Value = 62502250  # H83: 50002400 -> 62502250
$ws.Cells.Item(97, 8).Value = 1377.56  # H97: 1422.5416 -> 1377.56
$ws.Cells.Item(97, 9).Value = 878.0476  # I97: 907.05 -> 878.0476
$ws.Cells.Item(97, 11).Value = 878.0476  # K97: 907.05 -> 878.0476
$ws.Cells.Item(97, 13).Value = -382.0476  # M97: -411.05 -> -382.0476
$ws.Cells.Item(122, 8).Value = 4149.4165  # H122: 4232.75 -> 4149.4165
$ws.Cells.Item(122, 9).Value = 2764.2856  # I122: 2907.1428 -> 2764.2856
$ws.Cells.Item(122, 11).Value = 8292.856800000001  # K122: 8721.428400000001 -> 8292.856800000001
$ws.Cells.Item(122, 13).Value = -5842.856800000001  # M122: -6271.428400000001 -> -5842.856800000001
$ws.Cells.Item(123, 8).Value = 15988.75  # H123: 15524.533 -> 15988.75
$ws.Cells.Item(123, 10).Value = 15988.75  # J123: 15524.533 -> 15988.75
$ws.Cells.Item(123, 12).Value = 15988.75  # L123: 15524.533 -> 15988.75
$ws.Cells.Item(123, 14).Value = -20888.75  # N123: -20424.533 -> -20888.75
$ws.Cells.Item(124, 8).Value = 43780  # H124: 42780 -> 43780
$ws.Cells.Item(124, 10).Value = 43780  # J124: 42780 -> 43780
$ws.Cells.Item(124, 12).Value = 43780  # L124: 42780 -> 43780
$ws.Cells.Item(124, 14).Value = -53600  # N124: -52600 -> -53600
$ws.Cells.Item(126, 8).Value = 4100  # H126: 3996.44 -> 4100
$ws.Cells.Item(126, 9).Value = 2891.8918  # I126: 2881.7026 -> 2891.8918
$ws.Cells.Item(126, 10).Value = 5755.5557  # J126: 5081.8423 -> 5755.5557
$ws.Cells.Item(126, 11).Value = 8675.6754  # K126: 8645.1078 -> 8675.6754
$ws.Cells.Item(126, 12).Value = 17266.6671  # L126: 15245.5269 -> 17266.6671
$ws.Cells.Item(126, 13).Value = -6205.6754  # M126: -6175.1078 -> -6205.6754
$ws.Cells.Item(126, 14).Value = -22206.6671  # N126: -20185.5269 -> -22206.6671
$ws.Cells.Item(128, 8).Value = 42780  # H128: 41585 -> 42780
$ws.Cells.Item(128, 10).Value = 42780  # J128: 41585 -> 42780
$ws.Cells.Item(128, 12).Value = 42780  # L128: 41585 -> 42780
$ws.Cells.Item(128, 14).Value = -52740  # N128: -51545 -> -52740

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3122.487  # H122: 3147.9473 -> 3122.487
$ws.Cells.Item(122, 9).Value = 2527.9143  # I122: 2538.8823 -> 2527.9143
$ws.Cells.Item(122, 11).Value = 7583.742899999999  # K122: 7616.646900000001 -> 7583.742899999999
$ws.Cells.Item(122, 13).Value = -5133.742899999999  # M122: -5166.646900000001 -> -5133.742899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 133800  # H46: 117500 -> 133800
$ws.Cells.Item(46, 10).Value = 133800  # J46: 117500 -> 133800
$ws.Cells.Item(46, 12).Value = 133800  # L46: 117500 -> 133800
$ws.Cells.Item(46, 14).Value = -134262  # N46: -117962 -> -134262
$ws.Cells.Item(122, 8).Value = 4418.421  # H122: 5092.3335 -> 4418.421
$ws.Cells.Item(122, 9).Value = 2710.3572  # I122: 3512.077 -> 2710.3572
$ws.Cells.Item(122, 11).Value = 8131.071599999999  # K122: 10536.231 -> 8131.071599999999
$ws.Cells.Item(122, 13).Value = -5681.071599999999  # M122: -8086.231 -> -5681.071599999999
$ws.Cells.Item(126, 8).Value = 763859.1  # H126: 971658.6 -> 763859.1
$ws.Cells.Item(126, 9).Value = 2145.1667  # I126: 2497 -> 2145.1667
$ws.Cells.Item(126, 10).Value = 1335144.6  # J126: 1525465.2 -> 1335144.6
$ws.Cells.Item(126, 11).Value = 6435.500100000001  # K126: 7491 -> 6435.500100000001
$ws.Cells.Item(126, 12).Value = 4005433.8  # L126: 4576395.6 -> 4005433.8
$ws.Cells.Item(126, 13).Value = -3965.500100000001  # M126: -5021 -> -3965.500100000001
$ws.Cells.Item(126, 14).Value = -4010373.8  # N126: -4581335.6 -> -4010373.8
$ws.Cells.Item(134, 8).Value = 133800  # H134: 117500 -> 133800
$ws.Cells.Item(134, 10).Value = 133800  # J134: 117500 -> 133800
$ws.Cells.Item(134, 12).Value = 401400  # L134: 352500 -> 401400
$ws.Cells.Item(134, 14).Value = -406470  # N134: -357570 -> -406470
